# Auto-generated from the author diff: update crypto price/volume
# table cells (and the ImmutableX/NEARProtocol row swap) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.186.79"
$ws.Range("E2").Value = "'  +0.24%  "
$ws.Range("D3").Value = "'3.758.94"
$ws.Range("E3").Value = "'  +0.49%  "
$ws.Range("D5").Value = "'602.37"
$ws.Range("E5").Value = "'  +0.09%  "
$ws.Range("D6").Value = "'167.37"
$ws.Range("E6").Value = "'  -0.33%  "
$ws.Range("D7").Value = "'3.758.33"
$ws.Range("E7").Value = "'  +0.48%  "
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E9").Value = "'  +1.25%  "
$ws.Range("E10").Value = "'  +2.70%  "
$ws.Range("E11").Value = "'  +1.62%  "
$ws.Range("E12").Value = "'  -0.04%  "
$ws.Range("D13").Value = "'38.06"
$ws.Range("E13").Value = "'  -0.68%  "
$ws.Range("E14").Value = "'  +1.60%  "
$ws.Range("D15").Value = "'4.387.36"
$ws.Range("E15").Value = "'  +0.52%  "
$ws.Range("D16").Value = "'3.757.92"
$ws.Range("E16").Value = "'  +0.61%  "
$ws.Range("D17").Value = "'69.189.68"
$ws.Range("E17").Value = "'  +0.38%  "
$ws.Range("E18").Value = "'  +1.48%  "
$ws.Range("E19").Value = "'  +0.27%  "
$ws.Range("E20").Value = "'  -1.54%  "
$ws.Range("D21").Value = "'11.10"
$ws.Range("E21").Value = "'  +10.53%  "
$ws.Range("D22").Value = "'493.62"
$ws.Range("E22").Value = "'  -0.90%  "
$ws.Range("E23").Value = "'  +0.73%  "
$ws.Range("E24").Value = "'  +7.19%  "
$ws.Range("D25").Value = "'84.97"
$ws.Range("E25").Value = "'  +0.03%  "
$ws.Range("E26").Value = "'  -0.06%  "
$ws.Range("E27").Value = "'  +0.36%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "'  -0.14%  "
$ws.Range("E29").Value = "'  -0.10%  "
$ws.Range("E30").Value = "'  +1.25%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'8.16"
$ws.Range("E31").Value = "'  +2.66%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'2.48"
$ws.Range("E32").Value = "'  +2.27%  "
$ws.Range("D33").Value = "'31.62"
$ws.Range("E33").Value = "'  -0.48%  "
$ws.Range("D34").Value = "'3.904.49"
$ws.Range("E34").Value = "'  +0.63%  "
$ws.Range("D35").Value = "'3.691.19"
$ws.Range("E35").Value = "'  +0.59%  "
$ws.Range("E36").Value = "'  -0.37%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "'  +0.02%  "
$ws.Range("E38").Value = "'  +3.58%  "
$ws.Range("E39").Value = "'  +0.34%  "
$ws.Range("D40").Value = "'0.138"
$ws.Range("E40").Value = "'  +3.02%  "
$ws.Range("E41").Value = "'  +0.83%  "
$ws.Range("E42").Value = "'  +5.63%  "
$ws.Range("D43").Value = "'429.97"
$ws.Range("E43").Value = "'  -1.74%  "
$ws.Range("E44").Value = "'  -0.72%  "
$ws.Range("E45").Value = "'  -0.12%  "
$ws.Range("D46").Value = "'8.48"
$ws.Range("E46").Value = "'  +1.00%  "
$ws.Range("E48").Value = "'  -0.28%  "
$ws.Range("D49").Value = "'141.27"
$ws.Range("E49").Value = "'  -1.34%  "
$ws.Range("D50").Value = "'2.793.88"
$ws.Range("E50").Value = "'  +1.81%  "
$ws.Range("E51").Value = "'  +0.32%  "
